$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated extrapolation calibration values (removing < $5 price noise).
# Row => D, E, F, G, H values
$updates = @{
    3  = @(121091.9160649602, -0.04055452654159507, 0.2279851887159265, -3.576656112311131, 48.19106477650713)
    5  = @(122559.2999124138, -0.0390492214106278,  0.2275066541398969, -1.342150593475142, 13.41794205435837)
    6  = @(123285.471987498,  -0.03880019969021099, 0.2200654065035844, -0.7225663283659356, 8.24218565938005)
    7  = @(123823.5801999686, -0.04756574711378322, 0.2617606151182995, -1.834175766388056, 17.44050336386212)
    8  = @(125501.7654631743, -0.05542935253461787, 0.2116100857669064, -0.9484414421336175, 7.350683290569751)
    9  = @(127236.2742784592, -0.09643346658496654, 0.3298463083335351, -1.080162824515163, 8.125204318351777)
    10 = @(128589.3559845158, -0.1360349085765478,  0.4464757812208354, -1.6996865055798,   9.091009351829779)
    11 = @(130975.2902512479, -0.1462996400284794,  0.444399876863888,  -1.450693774417547, 7.245399232316193)
    14 = @(120257.1792292386, -0.07381262719815201, 0.193971534143024,  -0.6144705762785001, 8.433536171740222)
    16 = @(120259.0067400903, -0.1070763021382608,  0.2221697773350063, -0.3753743299946656, 6.142449630627786)
    17 = @(120223.878746479,  -0.09882300097050523, 0.2317129005192542, -0.6467046128420481, 8.771461899378204)
    18 = @(121122.606038857,  -0.0275323233216507,  0.1828769957451155, -0.7543573897679601, 6.811621619717825)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
